# Daily attendance processing - reorders the "Recorded By" (column G) list
# of contributors so the most recently-recorded-by entry is listed first
# (i.e. reverses the comma-separated list of names/emails in column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    if ($value -notlike "*,*") {
        continue
    }

    # Leave this specific combination untouched.
    if ($value -eq "System, admin@admin.com") {
        continue
    }

    $parts = $value.Split(",")
    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    $count = $trimmedParts.Length
    $reversedParts = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversedParts += $trimmedParts[$i]
    }

    $newValue = [string]::Join(", ", $reversedParts)
    $cell.Value2 = $newValue
}
